$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inventory_1.csv")

# Add new header "Machine" in N1 (match formatting of neighboring header cell M1)
$ws.Cells.Item(1, 14).Value = "Machine"
$ws.Range("N1").Font.Bold = $ws.Range("M1").Font.Bold
$ws.Range("N1").Font.Size = $ws.Range("M1").Font.Size
$ws.Range("N1").Borders.Item(9).LineStyle = $ws.Range("M1").Borders.Item(9).LineStyle

# Fill N2:N17 with 1
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 14).Value = 1
}

# Update selection to N18
$ws.Range("N18").Select()
